# Natmi following Dr Hou advice
# Rewrites the single data row (row 2) with corrected values and adds three
# new data rows (3-5) covering the additional "ECs" cluster combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("M2","Ccl12","Ccr1","M2",1,0.3333333333333333,3.288126333333333,9.864379,0.05813306630866938,0.05813306630866937,2,0.6666666666666666,114.5606336666667,343.681901,0.7368570786832789,0.736857078683279,376.6898363227199,3390.208526904479,0.04283576141510746,0.04283576141510746),
  @("M2","Ccl12","Ccr1","ECs",1,0.3333333333333333,3.288126333333333,9.864379,0.05813306630866938,0.05813306630866937,3,1,40.91135266666667,122.734058,0.263142921316721,0.2631429213167211,134.5216960355535,1210.695264319982,0.01529730489356191,0.01529730489356191),
  @("ECs","Ccl12","Ccr1","M2",3,1,53.27393966666666,159.821819,0.9418669336913307,0.9418669336913306,2,0.6666666666666666,114.5606336666667,343.681901,0.7368570786832789,0.736857078683279,6103.096286133102,54927.86657519792,0.6940213172681715,0.6940213172681715),
  @("ECs","Ccl12","Ccr1","ECs",3,1,53.27393966666666,159.821819,0.9418669336913307,0.9418669336913306,3,1,40.91135266666667,122.734058,0.263142921316721,0.2631429213167211,2179.508933645722,19615.5804028115,0.2478456164231591,0.2478456164231591)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $row = $rows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
